$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.861.23'
$ws.Range("E2").Value = '  +9.49%  '

$ws.Range("D3").Value = '3.463.67'
$ws.Range("E3").Value = '  +5.98%  '

$ws.Range("E4").Value = '  +0.11%  '

$ws.Range("D5").Value = '414.04'
$ws.Range("E5").Value = '  +3.84%  '

$ws.Range("D6").Value = '123.63'
$ws.Range("E6").Value = '  +13.44%  '

$ws.Range("D7").Value = '3.456.11'
$ws.Range("E7").Value = '  +5.91%  '

$ws.Range("D8").Value = '0.593'
$ws.Range("E8").Value = '  +2.29%  '

$ws.Range("E9").Value = '  +0.10%  '

$ws.Range("D10").Value = '0.682'
$ws.Range("E10").Value = '  +10.02%  '

$ws.Range("D11").Value = '0.129'
$ws.Range("E11").Value = '  +33.67%  '

$ws.Range("D12").Value = '41.21'
$ws.Range("E12").Value = '  +5.01%  '

$ws.Range("E13").Value = '  +0.42%  '

$ws.Range("D14").Value = '4.011.69'
$ws.Range("E14").Value = '  +6.06%  '

$ws.Range("E15").Value = '  +4.20%  '

$ws.Range("D16").Value = '19.91'
$ws.Range("E16").Value = '  +5.04%  '

$ws.Range("D17").Value = '3.466.76'
$ws.Range("E17").Value = '  +6.26%  '

$ws.Range("D18").Value = '62.791.88'
$ws.Range("E18").Value = '  +9.68%  '

$ws.Range("E19").Value = '  -0.29%  '

$ws.Range("D20").Value = '10.85'
$ws.Range("E20").Value = '  -2.38%  '

$ws.Range("D21").Value = '0.0000136'
$ws.Range("E21").Value = '  +26.60%  '

$ws.Range("D22").Value = '3.32'
$ws.Range("E22").Value = '  +0.08%  '

$ws.Range("D23").Value = '315.73'
$ws.Range("E23").Value = '  +6.07%  '

$ws.Range("D24").Value = '81.46'
$ws.Range("E24").Value = '  +9.71%  '

$ws.Range("D25").Value = '12.85'
$ws.Range("E25").Value = '  -0.50%  '

$ws.Range("E26").Value = '  -0.50%  '

$ws.Range("D27").Value = '30.81'
$ws.Range("E27").Value = '  +9.45%  '

$ws.Range("D28").Value = '7.78'
$ws.Range("E28").Value = '  +4.66%  '

$ws.Range("E29").Value = '  -0.78%  '

$ws.Range("E30").Value = '  -1.74%  '

$ws.Range("D31").Value = '0.175'
$ws.Range("E31").Value = '  +3.29%  '

$ws.Range("E32").Value = '  +4.06%  '

$ws.Range("D33").Value = '2.60'
$ws.Range("E33").Value = '  +21.56%  '

$ws.Range("D34").Value = '11.62'
$ws.Range("E34").Value = '  +3.39%  '

$ws.Range("B35").Value = 'InjectiveProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D35").Value = '42.23'
$ws.Range("E35").Value = '  +4.67%  '

$ws.Range("B36").Value = 'Dai'
$ws.Range("C36").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").Value = '  +0.03%  '

$ws.Range("E37").Value = '  -1.52%  '

$ws.Range("D38").Value = '52.21'
$ws.Range("E38").Value = '  +1.50%  '

$ws.Range("D39").Value = '0.998'
$ws.Range("E39").Value = '  -0.06%  '

$ws.Range("D40").Value = '3.50'
$ws.Range("E40").Value = '  +0.21%  '

$ws.Range("E41").Value = '  -2.36%  '

$ws.Range("E42").Value = '  +6.25%  '

$ws.Range("E43").Value = '  +3.28%  '

$ws.Range("D44").Value = '135.97'
$ws.Range("E44").Value = '  -1.52%  '

$ws.Range("D45").Value = '0.282'
$ws.Range("E45").Value = '  -0.76%  '

$ws.Range("D46").Value = '16.84'
$ws.Range("E46").Value = '  -0.08%  '

$ws.Range("D47").Value = '3.89'
$ws.Range("E47").Value = '  -0.14%  '

$ws.Range("E48").Value = '  +2.15%  '

$ws.Range("D49").Value = '21.96'
$ws.Range("E49").Value = '  -2.22%  '

$ws.Range("D50").Value = '2.207.40'
$ws.Range("E50").Value = '  +2.32%  '

$ws.Range("D51").Value = '2.48'
$ws.Range("E51").Value = '  +0.19%  '
